$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.41603216931191866, 0.41461622916590102),
    @(-0.34808618811088898, 0.34260388052615198),
    @(-0.14792723104268291, 0.14711292649773711),
    @(-0.13511292656638574, 0.13439207645659224),
    @(-0.12839207669860908, 0.12695877938107447),
    @(-0.025968792176827904, 0.025960334884598257),
    @(-0.0059603351791110981, 0.0059542154165637129),
    @(0.01404578428875336, -0.01408767992484794),
    @(0.020087679676915826, -0.02015050436326149),
    @(0.02615050411661457, -0.026165679922073082),
    @(0.030665679680378588, -0.030706525279665442),
    @(-0.045396631347287109, 0.045152169682640686),
    @(-0.039152169932540559, 0.039085519960493365),
    @(-0.027085520231297622, 0.02705317428894638),
    @(-0.021053174540974773, 0.021027700956528506),
    @(-0.015027701209403777, 0.015004576578685924),
    @(-0.0090045768326989517, 0.0089999997358711781),
    @(-0.056530041294539757, 0.056490339511032772),
    @(-0.027097341575204936, 0.027014002139807314),
    @(-0.018014002378434313, 0.018004316922233699),
    @(-0.0090043171612137485, 0.0089999997607748128),
    @(-0.18567664972955455, 0.18426817394059114),
    @(-0.12566523366214888, 0.12425171757615416),
    @(-0.042131046038297093, 0.041999999632198914),
    @(-0.094989987452649416, 0.094739124063746516),
    @(-0.088739124312883888, 0.088416684469656559),
    @(-0.082416684720258981, 0.08131671405579155),
    @(-0.075316714312127608, 0.074554428279172136),
    @(-0.062554428559074893, 0.062176038043318727),
    @(-0.042176038351379841, 0.042021051982597868),
    @(-0.027021052277047986, 0.027000963830575486),
    @(-0.0060009641451728157, 0.0059999997344588607),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
}
